$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "In Translation"
#    (shows up in the per-language status cells on every sheet)
# ------------------------------------------------------------------

# "Overview" sheet: zh-cn / de-de status columns (E & F), data rows 2-3
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2:F3").Value = "In Translation"

# "zh-cn" sheet: Status column (C), data rows 2-3
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2:C3").Value = "In Translation"

# "de-de" sheet: Status column (C), data rows 2-3
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2:C3").Value = "In Translation"

# ------------------------------------------------------------------
# 2. Narrow the status columns (they shrank from ~17.22 to ~13.41
#    characters wide). Excel quantizes ColumnWidth to whole pixels,
#    so 12.5 is the input that lands closest to the target width.
# ------------------------------------------------------------------

$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
